$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.131.36'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.904.40'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.92'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4609'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3929'
$ws.Range('E8').Value = '  +1.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.83'
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07930'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9998'
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.21'
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('D13').Value = '1.860.32'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.071'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.757'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06945'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.25'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.11'
$ws.Range('E20').Value = '  +1.86%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = '29.146.47'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.355'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '2.132.46'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.054'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.45'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.43'
$ws.Range('E28').Value = '  +1.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.093'
$ws.Range('E29').Value = '  +4.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.992'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '118.82'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09377'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9238'
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.325'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.269'
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.198'
$ws.Range('E37').Value = '  +3.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05819'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02103'
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.928'
$ws.Range('E40').Value = '  +3.18%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5741'
$ws.Range('E42').Value = '  +1.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1797'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.913'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.99'
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.234'
$ws.Range('E46').Value = '  +5.09%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5413'
$ws.Range('E47').Value = '  +2.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.07070'
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.876'
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.549'
$ws.Range('E50').Value = '  +5.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.10'
$ws.Range('E51').Value = '  -0.31%  '
